# Apply updated dSF (column F) values for the listed rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -3
    7  = -2
    8  = 4
    9  = -3
    10 = -3
    12 = -5
    13 = 3
    14 = 0
    17 = -7
    19 = 2
    20 = 7
    21 = -3
    26 = -1
    27 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
